# Scheduled-runner market data refresh: recompute leve profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for the affected leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 519.5102000000001
$ws.Range("J17").Value = 519.5102000000001
$ws.Range("L17").Value = 1558.5306
$ws.Range("N17").Value = -1894.5306
# Row 19
$ws.Range("H19").Value = 832
$ws.Range("J19").Value = 950.25
$ws.Range("L19").Value = 950.25
$ws.Range("N19").Value = -1300.25
# Row 132
$ws.Range("H132").Value = 779440.4
$ws.Range("I132").Value = 1407.6482
$ws.Range("J132").Value = 5447636.5
$ws.Range("K132").Value = 4222.944600000001
$ws.Range("L132").Value = 16342909.5
$ws.Range("M132").Value = -1692.944600000001
$ws.Range("N132").Value = -16347969.5
# Row 133
$ws.Range("H133").Value = 30000
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -40120
# Row 137
$ws.Range("H137").Value = 2042225.1
$ws.Range("I137").Value = 3126077.5
$ws.Range("J137").Value = 2032.3529
$ws.Range("K137").Value = 9378232.5
$ws.Range("L137").Value = 6097.0587
$ws.Range("M137").Value = -9375682.5
$ws.Range("N137").Value = -11197.0587
# Row 138
$ws.Range("H138").Value = 2875262.8
$ws.Range("I138").Value = 1215.8572
$ws.Range("J138").Value = 7248812.5
$ws.Range("K138").Value = 3647.5716
$ws.Range("L138").Value = 21746437.5
$ws.Range("M138").Value = 1492.4284
$ws.Range("N138").Value = -21756717.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1299.4706
$ws.Range("I45").Value = 930.0833
$ws.Range("J45").Value = 2186
$ws.Range("K45").Value = 930.0833
$ws.Range("L45").Value = 2186
$ws.Range("M45").Value = -553.0833
$ws.Range("N45").Value = -2940
# Row 61
$ws.Range("H61").Value = 43566200
$ws.Range("I61").Value = 66734056
$ws.Range("J61").Value = 126476.75
$ws.Range("K61").Value = 66734056
$ws.Range("L61").Value = 126476.75
$ws.Range("M61").Value = -66733844
$ws.Range("N61").Value = -126900.75
# Row 74
$ws.Range("H74").Value = 3938591.5
$ws.Range("I74").Value = 4922483
$ws.Range("J74").Value = 78710.62
$ws.Range("K74").Value = 4922483
$ws.Range("L74").Value = 78710.62
$ws.Range("M74").Value = -4921609
$ws.Range("N74").Value = -80458.62
# Row 77
$ws.Range("H77").Value = 3938591.5
$ws.Range("I77").Value = 4922483
$ws.Range("J77").Value = 78710.62
$ws.Range("K77").Value = 24612415
$ws.Range("L77").Value = 393553.1
$ws.Range("M77").Value = -24608047
$ws.Range("N77").Value = -402289.1
# Row 110
$ws.Range("H110").Value = 527580.4
$ws.Range("I110").Value = 715359.4
$ws.Range("J110").Value = 1799.2
$ws.Range("K110").Value = 715359.4
$ws.Range("L110").Value = 1799.2
$ws.Range("M110").Value = -713314.4
$ws.Range("N110").Value = -5889.2
# Row 118
$ws.Range("H118").Value = 60000
$ws.Range("J118").Value = 60000
$ws.Range("L118").Value = 60000
$ws.Range("N118").Value = -63314
# Row 132
$ws.Range("H132").Value = 68216.836
$ws.Range("I132").Value = 40660.32
$ws.Range("J132").Value = 183035.67
$ws.Range("K132").Value = 121980.96
$ws.Range("L132").Value = 549107.01
$ws.Range("M132").Value = -119450.96
$ws.Range("N132").Value = -554167.01
# Row 133
$ws.Range("H133").Value = 40174
$ws.Range("J133").Value = 40174
$ws.Range("L133").Value = 40174
$ws.Range("N133").Value = -45234
# Row 135
$ws.Range("H135").Value = 32317.309
$ws.Range("J135").Value = 32317.309
$ws.Range("L135").Value = 32317.309
$ws.Range("N135").Value = -42457.309
# Row 136
$ws.Range("H136").Value = 43566200
$ws.Range("I136").Value = 66734056
$ws.Range("J136").Value = 126476.75
$ws.Range("K136").Value = 200202168
$ws.Range("L136").Value = 379430.25
$ws.Range("M136").Value = -200199618
$ws.Range("N136").Value = -384530.25

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1558.1091
$ws.Range("I134").Value = 762.3171
$ws.Range("J134").Value = 3888.6428
$ws.Range("K134").Value = 2286.9513
$ws.Range("L134").Value = 11665.9284
$ws.Range("M134").Value = 248.0487000000003
$ws.Range("N134").Value = -16735.9284

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 24392120
$ws.Range("I58").Value = 31251718
$ws.Range("J58").Value = 2434.7778
$ws.Range("K58").Value = 31251718
$ws.Range("L58").Value = 2434.7778
$ws.Range("M58").Value = -31251515
$ws.Range("N58").Value = -2840.7778
# Row 132
$ws.Range("H132").Value = 38657.02
$ws.Range("I132").Value = 23093.2
$ws.Range("J132").Value = 126203.5
$ws.Range("K132").Value = 69279.60000000001
$ws.Range("L132").Value = 378610.5
$ws.Range("M132").Value = -66749.60000000001
$ws.Range("N132").Value = -383670.5
# Row 134
$ws.Range("H134").Value = 26713.28
$ws.Range("I134").Value = 1669.0333
$ws.Range("J134").Value = 84507.69500000001
$ws.Range("K134").Value = 5007.0999
$ws.Range("L134").Value = 253523.085
$ws.Range("M134").Value = -2472.0999
$ws.Range("N134").Value = -258593.085
# Row 136
$ws.Range("H136").Value = 24392120
$ws.Range("I136").Value = 31251718
$ws.Range("J136").Value = 2434.7778
$ws.Range("K136").Value = 93755154
$ws.Range("L136").Value = 7304.3334
$ws.Range("M136").Value = -93752604
$ws.Range("N136").Value = -12404.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 405.75
$ws.Range("I5").Value = 247.56
$ws.Range("J5").Value = 970.7143
$ws.Range("K5").Value = 742.6800000000001
$ws.Range("L5").Value = 2912.1429
$ws.Range("M5").Value = -630.6800000000001
$ws.Range("N5").Value = -3136.1429
# Row 80
$ws.Range("H80").Value = 1872
$ws.Range("J80").Value = 2250
$ws.Range("L80").Value = 6750
$ws.Range("N80").Value = -8622
# Row 83
$ws.Range("H83").Value = 1872
$ws.Range("J83").Value = 2250
$ws.Range("L83").Value = 20250
$ws.Range("N83").Value = -29610
# Row 107
$ws.Range("H107").Value = 968.17645
$ws.Range("I107").Value = 1220
$ws.Range("J107").Value = 684.875
$ws.Range("K107").Value = 3660
$ws.Range("L107").Value = 2054.625
$ws.Range("M107").Value = -1740
$ws.Range("N107").Value = -5894.625
# Row 113
$ws.Range("H113").Value = 486.0465
$ws.Range("I113").Value = 438.8889
$ws.Range("J113").Value = 728.5714
$ws.Range("K113").Value = 1316.6667
$ws.Range("L113").Value = 2185.7142
$ws.Range("M113").Value = 853.3333
$ws.Range("N113").Value = -6525.7142
# Row 115
$ws.Range("H115").Value = 2691.8667
$ws.Range("I115").Value = 2645.6
$ws.Range("K115").Value = 7936.799999999999
$ws.Range("M115").Value = -6761.799999999999
# Row 131
$ws.Range("H131").Value = 905.6604
$ws.Range("I131").Value = 445.92307
$ws.Range("J131").Value = 1055.075
$ws.Range("K131").Value = 1337.76921
$ws.Range("L131").Value = 3165.225
$ws.Range("M131").Value = 3702.23079
$ws.Range("N131").Value = -13245.225
# Row 135
$ws.Range("H135").Value = 405.75
$ws.Range("I135").Value = 247.56
$ws.Range("J135").Value = 970.7143
$ws.Range("K135").Value = 2228.04
$ws.Range("L135").Value = 8736.4287
$ws.Range("M135").Value = 306.96
$ws.Range("N135").Value = -13806.4287
# Row 138
$ws.Range("H138").Value = 3575.625
$ws.Range("I138").Value = 3258.5715
$ws.Range("J138").Value = 3822.2222
$ws.Range("K138").Value = 9775.7145
$ws.Range("L138").Value = 11466.6666
$ws.Range("M138").Value = -4635.7145
$ws.Range("N138").Value = -21746.6666

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 62252.79
$ws.Range("I132").Value = 41583.16
$ws.Range("J132").Value = 126845.375
$ws.Range("K132").Value = 124749.48
$ws.Range("L132").Value = 380536.125
$ws.Range("M132").Value = -122219.48
$ws.Range("N132").Value = -385596.125

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1336.409
$ws.Range("I100").Value = 1231.6316
$ws.Range("K100").Value = 1231.6316
$ws.Range("M100").Value = -690.6315999999999
# Row 132
$ws.Range("H132").Value = 25681.334
$ws.Range("I132").Value = 10444.632
$ws.Range("J132").Value = 170430
$ws.Range("K132").Value = 31333.896
$ws.Range("L132").Value = 511290
$ws.Range("M132").Value = -28803.896
$ws.Range("N132").Value = -516350

$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
# Row 132
$ws.Range("H132").Value = 51592.15
$ws.Range("I132").Value = 39230.92
$ws.Range("J132").Value = 123012.555
$ws.Range("K132").Value = 117692.76
$ws.Range("L132").Value = 369037.665
$ws.Range("M132").Value = -115162.76
$ws.Range("N132").Value = -374097.665
# Row 136
$ws.Range("H136").Value = 49257.74
$ws.Range("I136").Value = 35442.242
$ws.Range("J136").Value = 80076.92
$ws.Range("K136").Value = 106326.726
$ws.Range("L136").Value = 240230.76
$ws.Range("M136").Value = -103776.726
$ws.Range("N136").Value = -245330.76
